# CIERRE DE 6 DE DIC 2021
# Populate the November credits sheet ("REMISIONES   NOVIEMBRE  2021 ")
# with the newly-closed remisiones (rows 27-28 and 36-55), fix the
# "x" placeholder client name to "ISRAEL LEDO", and move the sheet's
# on-screen selection to reflect where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES   NOVIEMBRE  2021 ")
$ws.Activate()

# --- Rows that already had a client (D27/D28 were placeholder "x") ---
$ws.Range("A27").Value = 44515
$ws.Range("D27").Value = "COMERCIO CENTRAL "
$ws.Range("E27").Value = 3711

$ws.Range("A28").Value = 44515
$ws.Range("D28").Value = "OBRADOR"
$ws.Range("E28").Value = 2005

# --- Newly closed remisiones, rows 36-55 ---
$ws.Range("A36").Value = 44523
$ws.Range("D36").Value = "OBRADOR"
$ws.Range("E36").Value = 1189.81

$ws.Range("A37").Value = 44524
$ws.Range("D37").Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Range("E37").Value = 59045
$ws.Range("F37").Value = 44532
$ws.Range("G37").Value = 59045

$ws.Range("A38").Value = 44524
$ws.Range("D38").Value = "OBRADOR"
$ws.Range("E38").Value = 3655

$ws.Range("A39").Value = 44524
$ws.Range("D39").Value = "COMERCIO CENTRAL "
$ws.Range("E39").Value = 13576

$ws.Range("A40").Value = 44525
$ws.Range("D40").Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Range("E40").Value = 6961

$ws.Range("A41").Value = 44525
$ws.Range("D41").Value = "OBRADOR"
$ws.Range("E41").Value = 15564

$ws.Range("A42").Value = 44525
$ws.Range("D42").Value = "OBRADOR"
$ws.Range("E42").Value = 2279

$ws.Range("A43").Value = 44526
$ws.Range("D43").Value = "COMERCIO CENTRAL "
$ws.Range("E43").Value = 8799

$ws.Range("A44").Value = 44526
$ws.Range("D44").Value = "OBRADOR"
$ws.Range("E44").Value = 3337

$ws.Range("A45").Value = 44526
$ws.Range("D45").Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Range("E45").Value = 757

# Rows 46-47: cancelled remisiones - text + highlighted red/bold font
$ws.Range("A46").Value = 44527
$ws.Range("D46").Value = "CANCELADA"
$ws.Range("D46").Font.Bold = $true
$ws.Range("D46").Font.Size = 12
$ws.Range("D46").Font.Color = 255
$ws.Range("E46").Value = 0

$ws.Range("A47").Value = 44527
$ws.Range("D47").Value = "CANCELADA"
$ws.Range("D47").Font.Bold = $true
$ws.Range("D47").Font.Size = 12
$ws.Range("D47").Font.Color = 255
$ws.Range("E47").Value = 0

$ws.Range("A48").Value = 44529
$ws.Range("D48").Value = "COMERCIO CENTRAL "
$ws.Range("E48").Value = 376

$ws.Range("A49").Value = 44529
$ws.Range("D49").Value = "COMERCIO CENTRAL "
$ws.Range("E49").Value = 21811

$ws.Range("A50").Value = 44529
$ws.Range("D50").Value = "COMERCIO CENTRAL "
$ws.Range("E50").Value = 126

$ws.Range("A51").Value = 44529
$ws.Range("D51").Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Range("E51").Value = 53508
$ws.Range("F51").Value = 44532
$ws.Range("G51").Value = 53508

$ws.Range("A52").Value = 44529
$ws.Range("D52").Value = "COMERCIO CENTRAL "
$ws.Range("E52").Value = 86291

$ws.Range("A53").Value = 44529
$ws.Range("D53").Value = "COMERCIO CENTRAL "
$ws.Range("E53").Value = 161750

$ws.Range("A54").Value = 44531
$ws.Range("D54").Value = "ABASTOS DE 4 CARNES 11 SUR "
$ws.Range("E54").Value = 1421

$ws.Range("A55").Value = 44532
$ws.Range("D55").Value = "ISRAEL LEDO"
$ws.Range("E55").Value = 3135

# --- Cosmetic: page setup touched to portrait, scroll position + selection ---
$ws.PageSetup.Orientation = 1

$win = $excel.ActiveWindow
$win.ScrollRow = 33
$win.ScrollColumn = 1
$ws.Range("G52").Select()

Write-Host "Applied CIERRE DE 6 DE DIC 2021 updates"
